$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("I14").Value = 1
$ws.Range("J14").Value = 0.001
$ws.Range("K14").Value = 473
$ws.Range("L14").Value = 0.000946
